$wb = $excel.ActiveWorkbook

# --- Update the "sets" sheet: set 3's home_points goes from 23 to 24 ---
$sets = $wb.Worksheets.Item("sets")
$sets.Cells.Item(4, 4).Value = 24

# --- Append a new rally row (row 87) to the "rallies" sheet ---
$rallies = $wb.Worksheets.Item("rallies")
$newRow = 87

$rallies.Cells.Item($newRow, 1).Value = 86      # rally_id
$rallies.Cells.Item($newRow, 2).Value = 1       # match_id
$rallies.Cells.Item($newRow, 3).Value = 3       # set_number
$rallies.Cells.Item($newRow, 4).Value = 24      # rally_no
$rallies.Cells.Item($newRow, 5).Value = "NOS"   # side
$rallies.Cells.Item($newRow, 6).Value = "'"     # position (empty text cell)
$rallies.Cells.Item($newRow, 6).ClearFormats()  # drop the quote-prefix formatting so it stays a plain blank text cell
$rallies.Cells.Item($newRow, 7).Value = 5       # player_number
$rallies.Cells.Item($newRow, 8).Value = "LOB"   # action
$rallies.Cells.Item($newRow, 9).Value = "PONTO" # result
$rallies.Cells.Item($newRow, 10).Value = "NOS"  # who_scored
$rallies.Cells.Item($newRow, 11).Value = 24     # score_home
$rallies.Cells.Item($newRow, 12).Value = 0      # score_away
$rallies.Cells.Item($newRow, 13).Value = "1 5 lob"   # raw_text
$rallies.Cells.Item($newRow, 14).Value = "FRENTE"    # position_zone
$rallies.Cells.Item($newRow, 15).Value = "FRENTE"    # pos_fb
$rallies.Cells.Item($newRow, 16).Value = "FRENTE"    # frente_fundo
